# Update quarterly balance sheet with the newest quarter figures
# (adds the new "1402-04-14" period column values / relabels headers).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Header / period labels (row 9)
$ws.Range("I9").Value = "1402-04-14 (9)"
$ws.Range("M9").Value = "1402-04-14 (2)"

# Updated figures for the latest column (M)
$ws.Range("M14").Value = 33558574
$ws.Range("M18").Value = 64652484
$ws.Range("M27").Value = 116069676
$ws.Range("M35").Value = 56308173
$ws.Range("M37").Value = 73590840
$ws.Range("M43").Value = 77168184
$ws.Range("M56").Value = 36689608
$ws.Range("M57").Value = 38901492
$ws.Range("M58").Value = 116069676
